# Add the three new worksheets required for the CheckoutStepOne page object
# tests: verifyCartItemsDisplayCorrectly, verifyCartItemsCanBeRemoved,
# verifyCheckoutInformation. Each is appended after the last existing sheet.

$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd {
    param($workbook, $name)
    $lastSheet = $workbook.Worksheets.Item($workbook.Worksheets.Count)
    $newSheet = $workbook.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet.Name = $name
    return $newSheet
}

# ---------------------------------------------------------------------------
# Sheet: verifyCartItemsDisplayCorrectly
# ---------------------------------------------------------------------------
$s1 = Add-SheetAtEnd $wb "verifyCartItemsDisplayCorrectly"

$s1.Range("A1").Value = "itemsToAdd"
$s1.Range("A2").Value = "Sauce Labs Backpack"
$s1.Range("A3").Value = "Sauce Labs Fleece Jacket, Sauce Labs Bolt T-Shirt, Test.allTheThings() T-Shirt (Red)"
$s1.Range("A4").Value = "Sauce Labs Backpack, Sauce Labs Bike Light, Sauce Labs Bolt T-Shirt, Sauce Labs Onesie, Test.allTheThings() T-Shirt (Red)"

$s1.Range("B1").Value = "itemCount"
$s1.Range("B2").Value = "'1"
$s1.Range("B3").Value = "'3"
$s1.Range("B4").Value = "'5"

$s1.Columns.Item(1).ColumnWidth = 108.16666666666667
$s1.Columns.Item(2).ColumnWidth = 8.666666666666666

$s1.PageSetup.Orientation = 1

$s1.Range("A1:B4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: verifyCartItemsCanBeRemoved
# ---------------------------------------------------------------------------
$s2 = Add-SheetAtEnd $wb "verifyCartItemsCanBeRemoved"

$s2.Range("A1").Value = "itemsToAdd"
$s2.Range("B1").Value = "itemsToRemove"
$s2.Range("C1").Value = "itemCount"

$s2.Range("A2").Value = "Sauce Labs Backpack"
$s2.Range("B2").Value = "Sauce Labs Backpack"
$s2.Range("C2").Value = "'0"

$s2.Range("A3").Value = "Sauce Labs Fleece Jacket, Sauce Labs Bolt T-Shirt, Test.allTheThings() T-Shirt (Red)"
$s2.Range("B3").Value = "Sauce Labs Bolt T-Shirt"
$s2.Range("C3").Value = "'2"

$s2.Range("A4").Value = "Sauce Labs Backpack, Sauce Labs Bike Light, Sauce Labs Bolt T-Shirt, Sauce Labs Onesie, Test.allTheThings() T-Shirt (Red)"
$s2.Range("B4").Value = "Sauce Labs Backpack, Sauce Labs Bolt T-Shirt, Sauce Labs Onesie, Test.allTheThings() T-Shirt (Red)"
$s2.Range("C4").Value = "'1"

$s2.Columns.Item(1).ColumnWidth = 108.16666666666667
$s2.Columns.Item(2).ColumnWidth = 88.0
$s2.Columns.Item(3).ColumnWidth = 9.666666666666666

$s2.Range("F11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: verifyCheckoutInformation
# ---------------------------------------------------------------------------
$s3 = Add-SheetAtEnd $wb "verifyCheckoutInformation"

$s3.Range("A1").Value = "firstName"
$s3.Range("B1").Value = "lastName"
$s3.Range("C1").Value = "zipCode"
$s3.Range("D1").Value = "alert"

$s3.Range("A2").Value = "John"
$s3.Range("B2").Value = "Smith"
$s3.Range("C2").Value = 12345

$s3.Range("B3").Value = "Smith"
$s3.Range("C3").Value = 12345
$s3.Range("D3").Value = "Error: First Name is required"

$s3.Range("A4").Value = "John"
$s3.Range("C4").Value = 12345
$s3.Range("D4").Value = "Error: Last Name is required"

$s3.Range("A5").Value = "John"
$s3.Range("B5").Value = "Smith"
$s3.Range("D5").Value = "Error: Postal Code is required"

$s3.Columns.Item(1).ColumnWidth = 9.0
$s3.Columns.Item(2).ColumnWidth = 8.666666666666666
$s3.Columns.Item(3).ColumnWidth = 7.333333333333333
$s3.Columns.Item(4).ColumnWidth = 26.5

$s3.Range("G9").Select() | Out-Null
$s3.Activate()
